# Refresh the cryptocurrency price/volume snapshot (Price = column D, Volume(1h) = column E).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" values (column D) look like plain decimal numbers (e.g. "1.001").
# A leading apostrophe forces Excel to keep them as literal text, just like the
# original report (matching values such as "24.658.57" that are already non-numeric text).
function Set-TextValue($range, $text) {
    if ($text -match '^[+-]?\d+(\.\d+)?$') {
        $range.Value = "'" + $text
    } else {
        $range.Value = $text
    }
}

Set-TextValue $ws.Range("D2") '24.658.57'
Set-TextValue $ws.Range("E2") '  +3.94%  '
Set-TextValue $ws.Range("D3") '1.699.42'
Set-TextValue $ws.Range("E3") '  +2.43%  '
Set-TextValue $ws.Range("D4") '1.001'
Set-TextValue $ws.Range("E4") '  -0.23%  '
Set-TextValue $ws.Range("D5") '317.61'
Set-TextValue $ws.Range("E5") '  +3.18%  '
Set-TextValue $ws.Range("E6") '  -0.15%  '
Set-TextValue $ws.Range("D7") '0.3967'
Set-TextValue $ws.Range("E7") '  +2.36%  '
Set-TextValue $ws.Range("E8") '  +2.41%  '
Set-TextValue $ws.Range("D9") '1.547'
Set-TextValue $ws.Range("E9") '  +10.33%  '
Set-TextValue $ws.Range("D10") '54.70'
Set-TextValue $ws.Range("E10") '  +10.90%  '
Set-TextValue $ws.Range("E11") '  -0.27%  '
Set-TextValue $ws.Range("D12") '0.08835'
Set-TextValue $ws.Range("E12") '  +2.40%  '
Set-TextValue $ws.Range("D13") '7.314'
Set-TextValue $ws.Range("E13") '  +8.41%  '
Set-TextValue $ws.Range("D14") '23.45'
Set-TextValue $ws.Range("E14") '  +3.45%  '
Set-TextValue $ws.Range("E15") '  +2.11%  '
Set-TextValue $ws.Range("D16") '7.651'
Set-TextValue $ws.Range("E16") '  +6.33%  '
Set-TextValue $ws.Range("D17") '1.700.78'
Set-TextValue $ws.Range("E17") '  +2.11%  '
Set-TextValue $ws.Range("D18") '101.40'
Set-TextValue $ws.Range("E18") '  +1.66%  '
Set-TextValue $ws.Range("D19") '0.07099'
Set-TextValue $ws.Range("E19") '  +4.95%  '
Set-TextValue $ws.Range("D20") '19.82'
Set-TextValue $ws.Range("E20") '  +4.58%  '
Set-TextValue $ws.Range("D21") '6.893'
Set-TextValue $ws.Range("E21") '  +3.84%  '
Set-TextValue $ws.Range("D22") '1.000'
Set-TextValue $ws.Range("E22") '  -0.14%  '
Set-TextValue $ws.Range("D23") '14.18'
Set-TextValue $ws.Range("E23") '  +2.94%  '
Set-TextValue $ws.Range("D24") '24.643.99'
Set-TextValue $ws.Range("E24") '  +3.96%  '
Set-TextValue $ws.Range("D25") '3.088'
Set-TextValue $ws.Range("E25") '  +11.66%  '
Set-TextValue $ws.Range("D26") '2.331'
Set-TextValue $ws.Range("E26") '  +0.75%  '
Set-TextValue $ws.Range("D27") '22.50'
Set-TextValue $ws.Range("E27") '  +3.74%  '
Set-TextValue $ws.Range("D28") '160.20'
Set-TextValue $ws.Range("E28") '  +2.22%  '
Set-TextValue $ws.Range("D29") '5.240'
Set-TextValue $ws.Range("E29") '  +1.35%  '
Set-TextValue $ws.Range("D30") '134.59'
Set-TextValue $ws.Range("E30") '  +4.18%  '
Set-TextValue $ws.Range("D31") '7.638'
Set-TextValue $ws.Range("E31") '  +18.04%  '
Set-TextValue $ws.Range("E32") '  -0.27%  '
Set-TextValue $ws.Range("D33") '1.888.52'
Set-TextValue $ws.Range("E33") '  +2.17%  '
Set-TextValue $ws.Range("D34") '7.536'
Set-TextValue $ws.Range("E34") '  +16.67%  '
Set-TextValue $ws.Range("D35") '0.08616'
Set-TextValue $ws.Range("E35") '  +0.28%  '
Set-TextValue $ws.Range("D36") '11.58'
Set-TextValue $ws.Range("E36") '  +12.73%  '
Set-TextValue $ws.Range("D37") '0.2760'
Set-TextValue $ws.Range("E37") '  +4.64%  '
Set-TextValue $ws.Range("D38") '1.951'
Set-TextValue $ws.Range("E38") '  -1.36%  '
Set-TextValue $ws.Range("E39") '  +3.64%  '
Set-TextValue $ws.Range("D40") '0.02797'
Set-TextValue $ws.Range("E40") '  +11.08%  '
Set-TextValue $ws.Range("D41") '0.09095'
Set-TextValue $ws.Range("E41") '  +3.82%  '
Set-TextValue $ws.Range("D42") '0.7798'
Set-TextValue $ws.Range("E42") '  +3.94%  '
Set-TextValue $ws.Range("D43") '1.465'
Set-TextValue $ws.Range("D44") '0.7287'
Set-TextValue $ws.Range("E44") '  +4.23%  '
Set-TextValue $ws.Range("D45") '15.69'
Set-TextValue $ws.Range("E45") '  +5.67%  '
Set-TextValue $ws.Range("D46") '2.540'
Set-TextValue $ws.Range("E46") '  +7.16%  '
Set-TextValue $ws.Range("D47") '4.229'
Set-TextValue $ws.Range("E47") '  +4.09%  '
Set-TextValue $ws.Range("D48") '1.378'
Set-TextValue $ws.Range("E48") '  +16.57%  '
Set-TextValue $ws.Range("D50") '141.94'
Set-TextValue $ws.Range("E50") '  +2.13%  '
Set-TextValue $ws.Range("D51") '0.08060'
Set-TextValue $ws.Range("E51") '  +4.38%  '
